$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in X1 ("Data de criação"); copy V1's formatting (bold/border/centered header
# style) directly onto X1 so it reuses the same style definition, then set its text.
$ws.Range("V1").Copy($ws.Range("X1"))
$ws.Range("X1").Value = "Data de criação"

# "Início dos relatórios" (V) and "Término dos relatórios" (W) are the same for every row
$startDate = "30/09/2024"
$endDate = "29/10/2024"

# Creation date (X) per row, matching the date embedded in the ad name (column A)
$creationDates = @{
    2 = "16/09/2024"
    3 = "01/10/2024"
    4 = "16/10/2024"
    5 = "16/09/2024"
    6 = "16/09/2024"
    7 = "16/09/2024"
    8 = "15/10/2024"
    9 = "16/09/2024"
}

for ($row = 2; $row -le 9; $row++) {
    $vCell = $ws.Range("V$row")
    $wCell = $ws.Range("W$row")
    $xCell = $ws.Range("X$row")

    # Force text formatting so the dd/mm/yyyy strings are not reinterpreted as dates
    $vCell.NumberFormat = "@"
    $wCell.NumberFormat = "@"
    $xCell.NumberFormat = "@"

    $ws.Cells.Item($row, 22).Value = $startDate            # column V
    $ws.Cells.Item($row, 23).Value = $endDate              # column W
    $ws.Cells.Item($row, 24).Value = $creationDates[$row]  # column X

    # Restore the plain (unstyled) look used by the rest of the data cells
    $vCell.Style = $ws.Range("A$row").Style
    $wCell.Style = $ws.Range("A$row").Style
    $xCell.Style = $ws.Range("A$row").Style
}

$wb.Save()
